# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Integral" design / "Red Violet" colour
#                             scheme (the presentation's slide master /
#                             slides theme)
#   ppt/theme/theme2.xml  -> "Office Theme" design / "Office" colour
#                             scheme (the notes master's theme)
#
# The authored edit swaps the two themes' contents, so the slide
# master ends up with the "Office"/"Office Theme" colours and the
# notes master ends up with the "Red Violet"/"Integral" colours. The
# font scheme and format scheme are already identical between the two
# themes, so only the 12-slot colour scheme actually needs to change.
#
# Helper: turn an RRGGBB hex string into the RGB() long
# (R + G*256 + B*65536) that ThemeColorScheme items expect.
function ConvertTo-RgbLong([string]$Hex) {
    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office" colour scheme (currently on theme2.xml), expressed in
# ThemeColorScheme.Item(index) order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1..accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$p = $ppt.ActivePresentation

# Apply the "Office" colour scheme to the presentation's theme
# (ppt/theme/theme1.xml), replacing the current "Red Violet" scheme.
$master = $p.SlideMaster
$masterScheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le $masterScheme.Count; $i++) {
    $masterScheme.Item($i).RGB = ConvertTo-RgbLong $officeColors[$i - 1]
}
